# F03 Froze Encoder 1234
# Rewrites the per-epoch accuracy values in column B (re-run results after
# the encoder froze) and refreshes the stale Python object repr strings in
# column A for the "DisplayOutputs" rows (102-118), then updates the
# worksheet selection to match the new run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New accuracy values for column B, keyed by row number.
$newValues = @{
    4 = 0.953125
    5 = 0.96875
    6 = 0.9375
    7 = 0.890625
    8 = 0.890625
    9 = 0.875
    10 = 0.875
    11 = 0.921875
    12 = 0.875
    13 = 0.828125
    14 = 0.828125
    15 = 0.84375
    16 = 0.84375
    17 = 0.875
    18 = 0.828125
    19 = 0.828125
    22 = 0.828125
    23 = 0.8125
    24 = 0.8125
    25 = 0.8125
    26 = 0.8125
    32 = 0.796875
    33 = 0.796875
    34 = 0.796875
    39 = 0.78125
    40 = 0.78125
    41 = 0.78125
    42 = 0.78125
    43 = 0.78125
    44 = 0.78125
    45 = 0.78125
    46 = 0.78125
    47 = 0.78125
    51 = 0.796875
    52 = 0.796875
    53 = 0.796875
    54 = 0.796875
    55 = 0.796875
    56 = 0.796875
    57 = 0.796875
    58 = 0.796875
    59 = 0.796875
    60 = 0.796875
    61 = 0.796875
    62 = 0.796875
    63 = 0.796875
    64 = 0.796875
    69 = 0.796875
    70 = 0.796875
    71 = 0.8125
    72 = 0.8125
    73 = 0.8125
    74 = 0.8125
    75 = 0.828125
    76 = 0.828125
    77 = 0.828125
    78 = 0.828125
    79 = 0.828125
    80 = 0.828125
    81 = 0.828125
    82 = 0.828125
    83 = 0.828125
    84 = 0.828125
    85 = 0.828125
    86 = 0.828125
    87 = 0.828125
    88 = 0.828125
    89 = 0.828125
    90 = 0.828125
    91 = 0.828125
    92 = 0.828125
    93 = 0.828125
    94 = 0.828125
    95 = 0.828125
    96 = 0.828125
    97 = 0.828125
    98 = 0.828125
    99 = 0.828125
    100 = 0.828125
    101 = 0.828125
    102 = 0.828125
    103 = 0.875
    104 = 0.765625
    105 = 0.84375
    106 = 0.765625
    107 = 0.734375
    108 = 0.796875
    109 = 0.734375
    110 = 0.71875
    111 = 0.796875
    112 = 0.8125
    113 = 0.796875
    114 = 0.84375
    116 = 0.765625
    117 = 0.9375
    118 = 0.7377049180327869
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 2).Value = $newValues[$row]
}

# The "DisplayOutputs" object repr in column A (rows 102-118) changed because
# the notebook was re-run in a new process (new memory address) after the
# encoder was frozen.
$oldRepr = "<__main__.DisplayOutputs object at 0x7f565065bb50>"
$newRepr = "<__main__.DisplayOutputs object at 0x7f1450603be0>"
for ($row = 102; $row -le 118; $row++) {
    $ws.Cells.Item($row, 1).Value = $newRepr
}

# Update the active selection/scroll position left behind by the edit.
$ws.Range("H120").Select()
